$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 45707
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat

$ws.Range("C11").Value = "Tinte"
$ws.Range("D11").Value = 200
$ws.Range("E11").Value = "viriginia"
$ws.Range("G11").Value = "Efectivo"
